$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update invoice amount fields - these are stored as text in the sheet,
# so a leading apostrophe is used to force text entry and avoid Excel
# reinterpreting the decimal strings as numbers.
$ws.Range("I2").Value = "'48529.29"
$ws.Range("J2").Value = "'4416.85"
$ws.Range("K2").Value = "'1021.25"
$ws.Range("L2").Value = "'43091.19"

# Update the PDF thumbnail filename reference
$ws.Range("T2").Value = "captive aire_1754689853322.pdf"
